$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 99.14059966666667
$ws.Range("H2").Value = 297.421799
$ws.Range("I2").Value = 0.3911422343348016
$ws.Range("J2").Value = 0.3911422343348016
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2936666666666667
$ws.Range("N2").Value = 0.881
$ws.Range("O2").Value = 0.009113820319201367
$ws.Range("P2").Value = 0.009113820319201367
$ws.Range("Q2").Value = 29.11428943544445
$ws.Range("R2").Value = 262.028604919
$ws.Range("S2").Value = 0.003564800042978337
$ws.Range("T2").Value = 0.003564800042978337

$ws.Range("G3").Value = 99.14059966666667
$ws.Range("H3").Value = 297.421799
$ws.Range("I3").Value = 0.3911422343348016
$ws.Range("J3").Value = 0.3911422343348016
$ws.Range("O3").Value = 0.870405726797791
$ws.Range("P3").Value = 0.870405726797791
$ws.Range("Q3").Value = 2780.529280665026
$ws.Range("R3").Value = 25024.76352598523
$ws.Range("S3").Value = 0.3404524407574949
$ws.Range("T3").Value = 0.3404524407574949

$ws.Range("G4").Value = 99.14059966666667
$ws.Range("H4").Value = 297.421799
$ws.Range("I4").Value = 0.3911422343348016
$ws.Range("J4").Value = 0.3911422343348016
$ws.Range("O4").Value = 0.1204804528830076
$ws.Range("P4").Value = 0.1204804528830076
$ws.Range("Q4").Value = 384.8773240744213
$ws.Range("R4").Value = 3463.895916669791
$ws.Range("S4").Value = 0.04712499353432836
$ws.Range("T4").Value = 0.04712499353432837

$ws.Range("I5").Value = 0.4928190063160421
$ws.Range("J5").Value = 0.4928190063160421
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2936666666666667
$ws.Range("N5").Value = 0.881
$ws.Range("O5").Value = 0.009113820319201367
$ws.Range("P5").Value = 0.009113820319201367
$ws.Range("Q5").Value = 36.68250045555556
$ws.Range("R5").Value = 330.1425041
$ws.Range("S5").Value = 0.004491463873451771
$ws.Range("T5").Value = 0.004491463873451771

$ws.Range("I6").Value = 0.4928190063160421
$ws.Range("J6").Value = 0.4928190063160421
$ws.Range("O6").Value = 0.870405726797791
$ws.Range("P6").Value = 0.870405726797791
$ws.Range("S6").Value = 0.4289524853722798
$ws.Range("T6").Value = 0.4289524853722798

$ws.Range("I7").Value = 0.4928190063160421
$ws.Range("J7").Value = 0.4928190063160421
$ws.Range("O7").Value = 0.1204804528830076
$ws.Range("P7").Value = 0.1204804528830076
$ws.Range("S7").Value = 0.05937505707031052
$ws.Range("T7").Value = 0.05937505707031052

$ws.Range("G8").Value = 29.411685
$ws.Range("H8").Value = 88.235055
$ws.Range("I8").Value = 0.1160387593491562
$ws.Range("J8").Value = 0.1160387593491562
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2936666666666667
$ws.Range("N8").Value = 0.881
$ws.Range("O8").Value = 0.009113820319201367
$ws.Range("P8").Value = 0.009113820319201367
$ws.Range("Q8").Value = 8.637231495000002
$ws.Range("R8").Value = 77.73508345500001
$ws.Range("S8").Value = 0.001057556402771257
$ws.Range("T8").Value = 0.001057556402771257

$ws.Range("G9").Value = 29.411685
$ws.Range("H9").Value = 88.235055
$ws.Range("I9").Value = 0.1160387593491562
$ws.Range("J9").Value = 0.1160387593491562
$ws.Range("O9").Value = 0.870405726797791
$ws.Range("P9").Value = 0.870405726797791
$ws.Range("Q9").Value = 824.889617484255
$ws.Range("R9").Value = 7424.006557358295
$ws.Range("S9").Value = 0.1010008006680163
$ws.Range("T9").Value = 0.1010008006680163

$ws.Range("G10").Value = 29.411685
$ws.Range("H10").Value = 88.235055
$ws.Range("I10").Value = 0.1160387593491562
$ws.Range("J10").Value = 0.1160387593491562
$ws.Range("O10").Value = 0.1204804528830076
$ws.Range("P10").Value = 0.1204804528830076
$ws.Range("S10").Value = 0.01398040227836867
$ws.Range("T10").Value = 0.01398040227836867
